$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.411.02'
$ws.Range('E2').Value = '  +1.39%  '

$ws.Range('D3').Value = '1.824.53'
$ws.Range('E3').Value = '  +2.06%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.34'
$ws.Range('E5').Value = '  -0.11%  '

$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5352'
$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4021'
$ws.Range('E8').Value = '  +6.84%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07664'
$ws.Range('E9').Value = '  +2.84%  '

$ws.Range('E10').Value = '  +0.39%  '

$ws.Range('E11').Value = '  +1.57%  '

$ws.Range('E12').Value = '  +3.57%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.648'
$ws.Range('E13').Value = '  +5.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.001'
$ws.Range('E14').Value = '  +0.07%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.91'
$ws.Range('E15').Value = '  +1.57%  '

$ws.Range('D16').Value = '1.828.36'
$ws.Range('E16').Value = '  +2.84%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001078'
$ws.Range('E17').Value = '  +2.18%  '

$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.62'
$ws.Range('E18').Value = '  +0.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06585'
$ws.Range('E19').Value = '  +2.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.66'
$ws.Range('E20').Value = '  +2.09%  '

$ws.Range('E21').Value = '  +0.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.070'
$ws.Range('E22').Value = '  +2.92%  '

$ws.Range('D23').Value = '28.422.40'
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.09'
$ws.Range('E24').Value = '  -1.20%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.223'
$ws.Range('E25').Value = '  +6.48%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.456'
$ws.Range('E26').Value = '  +7.58%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.53'
$ws.Range('E27').Value = '  +1.51%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.69'
$ws.Range('E28').Value = '  +2.07%  '

$ws.Range('D29').Value = '2.038.25'
$ws.Range('E29').Value = '  +2.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.80'
$ws.Range('E30').Value = '  +3.11%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1115'
$ws.Range('E31').Value = '  +6.28%  '

$ws.Range('E32').Value = '  +1.13%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.673'
$ws.Range('E33').Value = '  +2.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07375'
$ws.Range('E34').Value = '  +14.56%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.643'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2242'
$ws.Range('E36').Value = '  -0.69%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02341'
$ws.Range('E37').Value = '  +2.34%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.227'
$ws.Range('E38').Value = '  +4.25%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.851'
$ws.Range('E39').Value = '  +4.53%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6265'
$ws.Range('E40').Value = '  +1.75%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.31'
$ws.Range('E41').Value = '  +1.92%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.180'
$ws.Range('E42').Value = '  +0.73%  '

$ws.Range('E43').Value = '  +0.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.393'
$ws.Range('E44').Value = '  -3.68%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.47'
$ws.Range('E45').Value = '  +2.13%  '

$ws.Range('E46').Value = '  +0.80%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5842'
$ws.Range('E47').Value = '  +1.38%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.71'
$ws.Range('E48').Value = '  -1.85%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.004'
$ws.Range('E49').Value = '  +4.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.203'
$ws.Range('E50').Value = '  +0.98%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06887'
$ws.Range('E51').Value = '  +1.25%  '
